# Pvo.Calc.xlsx edit: rename the leading "Index" column to "i" and
# renumber its values from a 1-based row index to a 0-based one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PVO")

# 1) Header rename: A1 "Index" -> "i" (also updates the "testdata" table's
#    first column name since A1 is the table's header cell).
$ws.Range("A1").Value2 = "i"

# 2) Re-index the data rows: column A currently holds 1..502 for rows 2..503;
#    shift every value down by one so it holds 0..501 (0-based index).
$dataRange = $ws.Range("A2:A503")
$vals = $dataRange.Value2
$rowCount = $vals.GetLength(0)
for ($r = 1; $r -le $rowCount; $r++) {
    $vals[$r, 1] = $vals[$r, 1] - 1
}
$dataRange.Value2 = $vals

# 3) Column A is narrower now that the header is "i" instead of "Index".
$ws.Columns.Item(1).ColumnWidth = 3.166666666666667
